# Edit Sendcodes.xlsx per commit: update messages list with new/renamed strings

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wording of the existing "Movement info" message (B8)
$ws.Range("B8").Value = "Movement info from clients"

# Add the new rows of message descriptions (B9:B16).
# Order matters here: new shared-string table entries are created in the
# order cells are first populated, so write them in the same sequence
# the saved workbook uses internally to keep the shared string table
# layout identical.
$ws.Range("B13").Value = "Send a message to host that the client is still in the lobby"
$ws.Range("B15").Value = "Add a new message to the chat box (client side)"
$ws.Range("B12").Value = "Update player locations (client side)"
$ws.Range("B11").Value = "Server gets a chat message and sends it to all of the players"
$ws.Range("B9").Value  = "Server recieves notice that the player is still in the lobby"
$ws.Range("B16").Value = "Setup game for client"
$ws.Range("B14").Value = "Delete a player that has left"
$ws.Range("B10").Value = "Player has told server that he is leaving, tell rest of players"

# Update the active selection to D1, matching the saved view state
$ws.Range("D1").Select()

$wb.Save()
